$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "all": append the 2020-05-10 (serial 43961) daily row as new row 33,
# pushing the footnote row down to row 34.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Rows.Item(33).Insert()
$wsAll.Range("A33").Value = 43961
$wsAll.Range("B33").Value = 277
$wsAll.Range("C33").Value = 272
$wsAll.Range("D33").Value = 77
$wsAll.Range("E33").Value = 67
$wsAll.Range("F33").Value = 10
$wsAll.Range("G33").Value = 8
$wsAll.Range("H33").Value = 187
$wsAll.Range("H34").Select()

# ---------------------------------------------------------------------------
# Sheet "kobe": correct the previous day (row 87) D/E values, then append
# the new day as row 88, pushing the footnote row down to row 89.
# ---------------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Range("D87").Value = 1
$wsKobe.Range("E87").Value = 277

$wsKobe.Rows.Item(88).Insert()
$wsKobe.Range("A88").Value = 43961
$wsKobe.Range("B88").Value = 0
$wsKobe.Range("C88").Value = 2562
$wsKobe.Range("D88").Value = 0
$wsKobe.Range("E88").Value = 277
$wsKobe.Range("F88").Value = 72
$wsKobe.Range("G88").Value = 63
$wsKobe.Range("H88").Value = 9
$wsKobe.Range("I88").Value = 8
$wsKobe.Range("J88").Value = 178
$wsKobe.Range("K88").Select()

# ---------------------------------------------------------------------------
# Sheet "other": append the new day as row 63, pushing the footnote row
# down to row 64. This sheet is the one that should end up active/selected,
# matching the original workbook's active tab.
# ---------------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Rows.Item(63).Insert()
$wsOther.Range("A63").Value = 43961
$wsOther.Range("B63").Value = 0
$wsOther.Range("C63").Value = 14
$wsOther.Range("D63").Value = 5
$wsOther.Range("E63").Value = 4
$wsOther.Range("F63").Value = 1
$wsOther.Range("G63").Value = 0
$wsOther.Range("H63").Value = 9
$wsOther.Range("A63").Select()

# Restore "other" as the active sheet (matches the workbook's original
# active tab / tabSelected state).
$wsOther.Activate()
